$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.057.25"
$ws.Range("E2").Value = "  -1.91%  "

$ws.Range("D3").Value = "1.832.15"
$ws.Range("E3").Value = "  -0.94%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.29%  "

$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4614"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3859"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07841"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9601"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.67%  "

$ws.Range("D12").Value = "1.863.78"
$ws.Range("E12").Value = "  +1.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.666"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.877"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06862"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009919"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "

$ws.Range("D21").Value = "28.086.63"
$ws.Range("E21").Value = "  -1.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.291"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.99%  "

$ws.Range("E23").Value = "  -2.78%  "

$ws.Range("E24").Value = "  -4.00%  "

$ws.Range("D25").Value = "2.093.94"
$ws.Range("E25").Value = "  +0.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.733"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.963"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9407"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09241"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.258"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.319"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.325"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05838"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02107"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.14%  "

$ws.Range("E38").Value = "  -1.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.704"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5588"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.885"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1757"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07299"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5267"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.137"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.109"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.830"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.022"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.07%  "
